# Update the Slovakia Covid Daily Stats workbook:
#  - Revise a handful of previously-published AgTests (H) / AgPosit (I)
#    values for existing rows (corrections from the new upstream file
#    structure effective 07.Feb.2021).
#  - Append a new day's record (row 337, date 2021-02-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Corrections to already-present rows (columns H = AgTests, I = AgPosit) ---

$ws.Range("H310").Value = 74794
$ws.Range("I310").Value = 3916

$ws.Range("H311").Value = 62461
$ws.Range("I311").Value = 1992

$ws.Range("I320").Value = 3696

$ws.Range("H323").Value = 149164

$ws.Range("H324").Value = 231289
$ws.Range("I324").Value = 2648

$ws.Range("H325").Value = 705128
$ws.Range("I325").Value = 5845

$ws.Range("H326").Value = 417163
$ws.Range("I326").Value = 3685

$ws.Range("H327").Value = 235537
$ws.Range("I327").Value = 2868

$ws.Range("H328").Value = 178276
$ws.Range("I328").Value = 2610

$ws.Range("H329").Value = 82192

$ws.Range("H331").Value = 147538
$ws.Range("I331").Value = 2558

$ws.Range("H332").Value = 411930
$ws.Range("I332").Value = 4026

$ws.Range("H333").Value = 252629
$ws.Range("I333").Value = 2705

$ws.Range("H334").Value = 201634
$ws.Range("I334").Value = 3349

$ws.Range("H335").Value = 121935
$ws.Range("I335").Value = 2798

$ws.Range("H336").Value = 96286
$ws.Range("I336").Value = 3105

# --- New row 337 (2021-02-04) ---

$ws.Range("A337").Value = 44231
$ws.Range("A337").NumberFormat = $ws.Range("A336").NumberFormat

$ws.Range("B337").Value = 259533
$ws.Range("C337").Value = 234371
$ws.Range("D337").Value = 20112
$ws.Range("E337").Value = 11282
$ws.Range("F337").Value = 2630
$ws.Range("G337").Value = 5050
$ws.Range("H337").Value = 93385
$ws.Range("I337").Value = 2979
